# Certificate of Employment template - add HTML frontend placeholders
# 1. Add a new blank centered/bold/sz36 paragraph at the very top of the body.
# 2. Introduce {pronoun_subject} and {pronoun_possessive} placeholders, replacing
#    the hard-coded "She"/"his" pronouns, and fix the missing space before
#    {start_date}.

$d = $word.ActiveDocument

# --- 1. Insert a brand-new empty paragraph before the current first paragraph ---
$firstPara = $d.Paragraphs.Item(1).Range
$firstPara.Collapse(1)            # wdCollapseStart
$firstPara.InsertParagraphBefore()

# --- 2. Fix "since{start_date}" -> "since {start_date}" (missing space) ---
$d.Content.Find.Execute("since{start_date}", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "since {start_date}", 2) | Out-Null

# --- 3. Replace the hard-coded "She" with the {pronoun_subject} placeholder ---
$d.Content.Find.Execute(". She is currently ", $true, $false, $false, $false, $false, `
                         $true, 1, $false, ". {pronoun_subject} is currently ", 2) | Out-Null

# --- 4. Replace the hard-coded "his" with the {pronoun_possessive} placeholder ---
$d.Content.Find.Execute("upon his request", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "upon {pronoun_possessive} request", 2) | Out-Null

Write-Host "Done. Paragraph count:" $d.Paragraphs.Count
